$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header in B1 from "Nombre humano" to "Sexo"
$ws.Range("B1").Value = "Sexo"

# Update the selected cell to match the saved view state
$ws.Range("E13").Select()
